# Documentation_Technique.docx edit script
# Commit: "Modification de l'introduction. Ajout de certaines informations sur le logiciel."
#
# 1) "...gerer leur billetterie de cinema." -> "...gerer leur propre billetterie de cinema."
# 2) "...les films, connaitre le..." -> "...les films, les seances et connaitre le..."
# 3) Split the introduction paragraph: remove the trailing space after
#    "nombre de place vendue." and start a new paragraph describing the
#    two modes (user / administrator).
# 4) Move the "_GoBack" bookmark from the old Statistique paragraph to the
#    end of the newly added paragraph (this is what Word does naturally
#    since _GoBack always tracks the very last edit location).
# 5) Refresh the cached "PAGE" field result shown in the footer.

$d = $word.ActiveDocument

# --- 1) "leur billetterie de cinema." -> "leur propre billetterie de cinema." ---
$rng = $d.Content
$null = $rng.Find.Execute(
    "leur billetterie de cinéma.", $true, $false, $false, $false, $false,
    $true, 1, $false, "leur propre billetterie de cinéma.", 2)

# --- 2) "les films, connaitre le" -> "les films, les seances et connaitre le" ---
$rng = $d.Content
$null = $rng.Find.Execute(
    "les films, connaitre le", $true, $false, $false, $false, $false,
    $true, 1, $false, "les films, les séances et connaitre le", 2)

# --- 3) Split the paragraph: drop the trailing space, start a new paragraph ---
$rng = $d.Content
$null = $rng.Find.Execute(
    "nombre de place vendue. ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nombre de place vendue.^pL’application est divisée en deux modes. Le mode utilisateur qui donne la possibilité de faire une réservation de séance et un mode administrateur qui lui permet de gérer les films, les salles et les séances. L’administrateur a un privilège plus grand que l’utilisateur et a donc accès à plus de fonctionnalités.",
    2)

# --- 4) Move the "_GoBack" bookmark to right before the final period of
#        the sentence we just typed (that is where editing stopped). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$rng = $d.Content
$null = $rng.Find.Execute(
    "plus de fonctionnalités.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$goBackPoint = $d.Range($rng.End - 1, $rng.End - 1)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# --- 5) Refresh the cached PAGE field result in the footer ---
$ftr = $d.Sections(1).Footers(1)
$frng = $ftr.Range
$null = $frng.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)

Write-Output "edit complete"
